# paises.xlsx COVID-19 data refresh: updates the 'Pais' worksheet from the
# 5 Jul 2020 14:10 snapshot to the 15:27 snapshot (updated case counts) and
# re-sorts a handful of rows whose totals changed enough to swap rank order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (A1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 5 de Julio de 2020 a las 15:27"

# Straight value refreshes (country stays in the same row)
# Row 4
$ws.Cells.Item(4, 2).Value = 2936890
$ws.Cells.Item(4, 3).Value = 1120
$ws.Cells.Item(4, 5).Value = 1543940
$ws.Cells.Item(4, 7).Value = 13
$ws.Cells.Item(4, 8).Value = 132331
# Row 16
$ws.Cells.Item(16, 2).Value = 209509
$ws.Cells.Item(16, 3).Value = 3580
$ws.Cells.Item(16, 4).Value = 145236
$ws.Cells.Item(16, 5).Value = 62357
$ws.Cells.Item(16, 7).Value = 58
$ws.Cells.Item(16, 8).Value = 1916
# Row 24
$ws.Cells.Item(24, 2).Value = 99799
$ws.Cells.Item(24, 3).Value = 616
$ws.Cells.Item(24, 4).Value = 92284
$ws.Cells.Item(24, 5).Value = 7387
$ws.Cells.Item(24, 7).Value = 5
$ws.Cells.Item(24, 8).Value = 128
# Row 26
$ws.Cells.Item(26, 4).Value = 27597
$ws.Cells.Item(26, 5).Value = 46289
$ws.Cells.Item(26, 7).Value = 9
$ws.Cells.Item(26, 8).Value = 1490
# Row 36
$ws.Cells.Item(36, 2).Value = 49941
$ws.Cells.Item(36, 3).Value = 638
$ws.Cells.Item(36, 4).Value = 40463
$ws.Cells.Item(36, 5).Value = 9110
$ws.Cells.Item(36, 7).Value = 3
$ws.Cells.Item(36, 8).Value = 368
# Row 97
$ws.Cells.Item(97, 5).Value = 3564
$ws.Cells.Item(97, 7).Value = 3
$ws.Cells.Item(97, 8).Value = 16
# Row 101
$ws.Cells.Item(101, 2).Value = 3151
$ws.Cells.Item(101, 3).Value = 57
$ws.Cells.Item(101, 4).Value = 2196
$ws.Cells.Item(101, 5).Value = 842
# Row 117
$ws.Cells.Item(117, 2).Value = 1836
$ws.Cells.Item(117, 3).Value = 5
$ws.Cells.Item(117, 5).Value = 212
# Row 127
$ws.Cells.Item(127, 2).Value = 1269
$ws.Cells.Item(127, 3).Value = 10
$ws.Cells.Item(127, 4).Value = 1156
$ws.Cells.Item(127, 5).Value = 106
# Row 132
$ws.Cells.Item(132, 2).Value = 1124
$ws.Cells.Item(132, 3).Value = 1
$ws.Cells.Item(132, 5).Value = 94
# Row 153
$ws.Cells.Item(153, 4).Value = 652
$ws.Cells.Item(153, 5).Value = 11
# Row 175
$ws.Cells.Item(175, 2).Value = 179
$ws.Cells.Item(175, 3).Value = 1
$ws.Cells.Item(175, 5).Value = 3

# Rank shuffles: totals moved enough that countries swap row order.
# Rows 76-80 (~7.3k-7.8k cases): El Salvador now leads, pushing
# Kenia/Senegal/Consejo Danes/Kirguistan down a slot each.
# Row 76
$ws.Cells.Item(76, 1).Value = "El Salvador"
$ws.Cells.Item(76, 2).Value = 7777
$ws.Cells.Item(76, 3).Value = 270
$ws.Cells.Item(76, 4).Value = 4573
$ws.Cells.Item(76, 5).Value = 2987
$ws.Cells.Item(76, 7).Value = 7
$ws.Cells.Item(76, 8).Value = 217
# Row 77
$ws.Cells.Item(77, 1).Value = "Kenia"
$ws.Cells.Item(77, 2).Value = 7577
$ws.Cells.Item(77, 4).Value = 2236
$ws.Cells.Item(77, 5).Value = 5182
$ws.Cells.Item(77, 8).Value = 159
# Row 78
$ws.Cells.Item(78, 1).Value = "Senegal"
$ws.Cells.Item(78, 2).Value = 7400
$ws.Cells.Item(78, 3).Value = 128
$ws.Cells.Item(78, 4).Value = 4870
$ws.Cells.Item(78, 5).Value = 2397
$ws.Cells.Item(78, 7).Value = 4
$ws.Cells.Item(78, 8).Value = 133
# Row 79
$ws.Cells.Item(79, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(79, 2).Value = 7379
$ws.Cells.Item(79, 3).Value = 0
$ws.Cells.Item(79, 4).Value = 2961
$ws.Cells.Item(79, 5).Value = 4236
$ws.Cells.Item(79, 7).Value = 0
$ws.Cells.Item(79, 8).Value = 182
# Row 80
$ws.Cells.Item(80, 1).Value = "Kirguistan"
$ws.Cells.Item(80, 2).Value = 7377
$ws.Cells.Item(80, 3).Value = 283
$ws.Cells.Item(80, 4).Value = 2802
$ws.Cells.Item(80, 5).Value = 4487
$ws.Cells.Item(80, 7).Value = 10
$ws.Cells.Item(80, 8).Value = 88

# Rows 159-162 (~340-360 cases): Siria moves ahead of Vietnam/Angola/Mauricio.
# Row 159
$ws.Cells.Item(159, 1).Value = "Siria"
$ws.Cells.Item(159, 2).Value = 358
$ws.Cells.Item(159, 3).Value = 20
$ws.Cells.Item(159, 4).Value = 126
$ws.Cells.Item(159, 5).Value = 219
$ws.Cells.Item(159, 7).Value = 3
$ws.Cells.Item(159, 8).Value = 13
# Row 160
$ws.Cells.Item(160, 1).Value = "Vietnam"
$ws.Cells.Item(160, 2).Value = 355
$ws.Cells.Item(160, 4).Value = 340
$ws.Cells.Item(160, 5).Value = 15
$ws.Cells.Item(160, 8).Value = 0
# Row 161
$ws.Cells.Item(161, 1).Value = "Angola"
$ws.Cells.Item(161, 2).Value = 346
$ws.Cells.Item(161, 4).Value = 108
$ws.Cells.Item(161, 5).Value = 219
$ws.Cells.Item(161, 8).Value = 19
# Row 162
$ws.Cells.Item(162, 1).Value = "Mauricio"
$ws.Cells.Item(162, 2).Value = 341
$ws.Cells.Item(162, 4).Value = 330
$ws.Cells.Item(162, 5).Value = 1

# Rows 205-206 (tied at 18 cases): Fiyi/Dominica swap.
$ws.Cells.Item(205, 1).Value = "Fiyi"
$ws.Cells.Item(206, 1).Value = "Dominica"
